# Applies the "minor changes to values" commit:
#   - Image Size values for the three EfficientNet rows (H7:H9) are corrected
#     from the stray "366x367" / "366x368" / "366x369" strings to the same
#     "366x366" value used elsewhere in the sheet. This also makes the three
#     now-unused strings drop out of the shared strings table, which shifts
#     the shared-string indices used by the metric headers in row 1
#     (J1:O1) down by three automatically when the workbook is saved.
#   - The active selection on the sheet moves to R13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = "366x366"
$ws.Range("H8").Value = "366x366"
$ws.Range("H9").Value = "366x366"

$ws.Range("R13").Select()
